$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Insert a new row at 14 (everything at/after row 14 shifts down by one). ---
# The new row inherits formatting from the row above (B13, style s="5"),
# matching rows 12/13/15/16 which are all the "note" style.
[void]$ws.Rows.Item(14).Insert()

# --- 2. Apply the formatting for the new "Total" row and "Score" column first
#        (so the new cell-style entries land in the same order as the reference
#        workbook); values are filled in afterwards. ---
$ws.Range("D19").HorizontalAlignment = -4152

[void]$ws.Range("B18").Copy()
[void]$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("E19").NumberFormat = "#,##0.0"

[void]$ws.Range("D20").Copy()
[void]$ws.Range("E20").PasteSpecial(-4122)

[void]$ws.Range("B21").Copy()
[void]$ws.Range("E21").PasteSpecial(-4122)
$ws.Range("E21").NumberFormat = "#,##0.0"
$ws.Range("E21").HorizontalAlignment = -4152

$ws.Application.CutCopyMode = $false

# --- 3. Fill in the cell contents. ---
# Edit the note text in B13 (now reworded) and fill the newly inserted B14.
$ws.Range("B13").Value = "* Do not use a formula in a cell that may have its position changed after the merge (for example under a TBS block). Otherwise Excel will raise an error message."
$ws.Range("B14").Value = "    This is because the location of formulas are saved a second time in another sub-file for the order of evaluation."
# B15 / B16 already hold the correct (shifted) content, no edit needed there.

# New "Score" column header + placeholder field, and the "Total:" / SUM formula row.
$ws.Range("E20").Value = "Score"
$ws.Range("E21").Value = "[a.score;ope=xlsxNum]"
$ws.Range("D19").Value = "Total:"
$ws.Range("E19").Formula = "=SUM(E21:E2000)"

# --- 4. Move the active selection to E20, like in the edited workbook. ---
[void]$ws.Range("E20").Select()
